# [MOD] Adjusted size for text fields and scrollpane
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the literal "null" placeholder text left over from the table import ---
$nullCells = @(
    "C2",
    "B3", "D3",
    "C4", "D4",
    "C5",
    "C6",
    "C7",
    "C8",
    "B9", "D9",
    "C10", "D10",
    "C11",
    "C12",
    "B13", "D13",
    "C15", "D15",
    "C16", "D16",
    "C17", "D17",
    "C18",
    "C19"
)
foreach ($addr in $nullCells) {
    $ws.Range($addr).Value = ""
}

# --- Row 5: Ctenophore -> Beroe forskalii, and tweak its description ---
$ws.Range("B5").Value = "Beroe forskalii"
$ws.Range("E5").Value = "Questa specie di ctenoforo di piccole dimensioni presenta un corpo completamente trasparente, fatta eccezione di filamenti “dentati” all’interno del corpo"

# --- Fill in the two previously-empty records (rows 20 and 21) ---
$ws.Range("B20").Value = "Lampocteis cruentiventer"
$ws.Range("D20").Value = "Ctenoforo dal ventre insanguinato"
$ws.Range("E20").Value = "Questo ctenoforo presenta un colore rosso molto intenso e dei filamenti interni seghettati che vengono talvolta attraversati da particelle bio-luminose sconosciute. Al momento sono l’unica specie conosciuta del loro genere."

$ws.Range("B21").Value = "Hexatrygon bickelli"
$ws.Range("D21").Value = "Razza esabranchiata"
$ws.Range("E21").Value = "Questa specie di razza presenta 6 paia di branchie, da cui deriva il suo nome. La testa è allungata e le dimensioni totali possono raggiungere anche 1.7m. Solitamente vivono a stretto contatto con il fondale marino."

# --- Row height adjustments (text fields resized) ---
$ws.Rows.Item(3).RowHeight = 27.7
$ws.Rows.Item(5).RowHeight = 27.7
$ws.Rows.Item(8).RowHeight = 27.7
$ws.Rows.Item(20).RowHeight = 40.95
$ws.Rows.Item(21).RowHeight = 40.95

# --- Scrollpane / selection moved ---
$ws.Range("E26").Select()
